$d = $word.ActiveDocument

$replacements = @(
    @{old = "99×99="; new = "22×47="},
    @{old = "84×17="; new = "92×27="},
    @{old = "39×98="; new = "74×49="},
    @{old = "43×69="; new = "16×61="},
    @{old = "18×60="; new = "88×42="},
    @{old = "49×78="; new = "88×41="},
    @{old = "99×52="; new = "98×37="},
    @{old = "54×57="; new = "24×11="},
    @{old = "83×20="; new = "19×95="},
    @{old = "61×36="; new = "56×71="},
    @{old = "55×44="; new = "64×38="},
    @{old = "11×52="; new = "63×96="},
    @{old = "46×48="; new = "35×48="},
    @{old = "78×14="; new = "78×62="},
    @{old = "95×22="; new = "59×62="},
    @{old = "87×40="; new = "66×27="},
    @{old = "99×36="; new = "23×18="},
    @{old = "30×84="; new = "41×84="},
    @{old = "33×47="; new = "31×12="},
    @{old = "82×17="; new = "15×31="},
    @{old = "16×43="; new = "39×90="},
    @{old = "99×71="; new = "68×25="},
    @{old = "77×28="; new = "11×16="},
    @{old = "97×20="; new = "78×29="},
    @{old = "97×86="; new = "14×57="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
